$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of row 5 into row 6
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = "2025-11-07 06:49:21 UTC"
$ws.Range("B6").Value = "2025-11-07 12:19:21 IST"
$ws.Range("C6").Value = "SKIPPED"
$ws.Range("D6").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E6").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = ""
